$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 216.11111
$ws.Cells.Item(33, 9).Value = 111.318184
$ws.Cells.Item(33, 10).Value = 677.2
$ws.Cells.Item(33, 11).Value = 111.318184
$ws.Cells.Item(33, 12).Value = 677.2
$ws.Cells.Item(33, 13).Value = 117.681816
$ws.Cells.Item(33, 14).Value = -1135.2

$ws.Cells.Item(43, 8).Value = 2454
$ws.Cells.Item(43, 9).Value = 2464
$ws.Cells.Item(43, 10).Value = 2446.5
$ws.Cells.Item(43, 11).Value = 2464
$ws.Cells.Item(43, 12).Value = 2446.5
$ws.Cells.Item(43, 13).Value = -2395

$ws.Cells.Item(64, 8).Value = 6796.6
$ws.Cells.Item(64, 9).Value = 7495
$ws.Cells.Item(64, 10).Value = 6331
$ws.Cells.Item(64, 11).Value = 7495
$ws.Cells.Item(64, 12).Value = 6331
$ws.Cells.Item(64, 13).Value = -7247

$ws.Cells.Item(67, 8).Value = 6796.6
$ws.Cells.Item(67, 9).Value = 7495
$ws.Cells.Item(67, 10).Value = 6331
$ws.Cells.Item(67, 11).Value = 7495
$ws.Cells.Item(67, 12).Value = 6331
$ws.Cells.Item(67, 13).Value = -6637

$ws.Cells.Item(70, 8).Value = 2867.8572
$ws.Cells.Item(70, 9).Value = 1300
$ws.Cells.Item(70, 10).Value = 3495
$ws.Cells.Item(70, 11).Value = 3900
$ws.Cells.Item(70, 12).Value = 10485
$ws.Cells.Item(70, 13).Value = -3630
$ws.Cells.Item(70, 14).Value = -11025

$ws.Cells.Item(73, 8).Value = 2867.8572
$ws.Cells.Item(73, 9).Value = 1300
$ws.Cells.Item(73, 10).Value = 3495
$ws.Cells.Item(73, 11).Value = 3900
$ws.Cells.Item(73, 12).Value = 10485
$ws.Cells.Item(73, 13).Value = -2964
$ws.Cells.Item(73, 14).Value = -12357

$ws.Cells.Item(74, 8).Value = 7453.875
$ws.Cells.Item(74, 9).Value = 5003
$ws.Cells.Item(74, 10).Value = 7804
$ws.Cells.Item(74, 11).Value = 5003
$ws.Cells.Item(74, 12).Value = 7804
$ws.Cells.Item(74, 13).Value = -4067
$ws.Cells.Item(74, 14).Value = -9676

$ws.Cells.Item(77, 8).Value = 7453.875
$ws.Cells.Item(77, 9).Value = 5003
$ws.Cells.Item(77, 10).Value = 7804
$ws.Cells.Item(77, 11).Value = 25015
$ws.Cells.Item(77, 12).Value = 39020
$ws.Cells.Item(77, 13).Value = -20335
$ws.Cells.Item(77, 14).Value = -48380

$ws.Cells.Item(106, 8).Value = 2488.3333
$ws.Cells.Item(106, 9).Value = 2532.5
$ws.Cells.Item(106, 10).Value = 2400
$ws.Cells.Item(106, 11).Value = 2532.5
$ws.Cells.Item(106, 12).Value = 2400
$ws.Cells.Item(106, 13).Value = -1901.5

$ws.Cells.Item(116, 8).Value = 378988.1
$ws.Cells.Item(116, 9).Value = 85345.14
$ws.Cells.Item(116, 10).Value = 892863.25
$ws.Cells.Item(116, 11).Value = 85345.14
$ws.Cells.Item(116, 12).Value = 892863.25
$ws.Cells.Item(116, 13).Value = -81903.14

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 17604.584
$ws.Cells.Item(32, 9).Value = 18445.371
$ws.Cells.Item(32, 10).Value = 228.33333
$ws.Cells.Item(32, 11).Value = 18445.371
$ws.Cells.Item(32, 12).Value = 228.33333
$ws.Cells.Item(32, 13).Value = -18158.371

$ws.Cells.Item(45, 8).Value = 3617.6924
$ws.Cells.Item(45, 9).Value = 2447.7778
$ws.Cells.Item(45, 10).Value = 6250
$ws.Cells.Item(45, 11).Value = 2447.7778
$ws.Cells.Item(45, 12).Value = 6250
$ws.Cells.Item(45, 13).Value = -2070.7778

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(55, 8).Value = 57000
$ws.Cells.Item(55, 9).Value = 0
$ws.Cells.Item(55, 10).Value = 57000
$ws.Cells.Item(55, 11).Value = 0
$ws.Cells.Item(55, 12).Value = 57000
$ws.Cells.Item(55, 14).Value = -57546

$ws.Cells.Item(80, 8).Value = 2396.6875
$ws.Cells.Item(80, 9).Value = 694
$ws.Cells.Item(80, 10).Value = 2639.9285
$ws.Cells.Item(80, 11).Value = 694
$ws.Cells.Item(80, 12).Value = 2639.9285
$ws.Cells.Item(80, 13).Value = 304
$ws.Cells.Item(80, 14).Value = -4635.9285

$ws.Cells.Item(83, 8).Value = 2396.6875
$ws.Cells.Item(83, 9).Value = 694
$ws.Cells.Item(83, 10).Value = 2639.9285
$ws.Cells.Item(83, 11).Value = 3470
$ws.Cells.Item(83, 12).Value = 13199.6425
$ws.Cells.Item(83, 13).Value = 1522
$ws.Cells.Item(83, 14).Value = -23183.6425

$ws.Cells.Item(134, 8).Value = 730184.0600000001
$ws.Cells.Item(134, 9).Value = 662992.25
$ws.Cells.Item(134, 10).Value = 1536485.4
$ws.Cells.Item(134, 11).Value = 1988976.75
$ws.Cells.Item(134, 12).Value = 4609456.199999999
$ws.Cells.Item(134, 13).Value = -1986441.75
$ws.Cells.Item(134, 14).Value = -4614526.199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 536.8333
$ws.Cells.Item(22, 9).Value = 616
$ws.Cells.Item(22, 10).Value = 497.25
$ws.Cells.Item(22, 11).Value = 616
$ws.Cells.Item(22, 12).Value = 497.25
$ws.Cells.Item(22, 13).Value = -266

$ws.Cells.Item(31, 8).Value = 10955.667
$ws.Cells.Item(31, 9).Value = 3857.7354
$ws.Cells.Item(31, 10).Value = 28193.5
$ws.Cells.Item(31, 11).Value = 3857.7354
$ws.Cells.Item(31, 12).Value = 28193.5
$ws.Cells.Item(31, 13).Value = -3562.7354

$ws.Cells.Item(34, 8).Value = 10955.667
$ws.Cells.Item(34, 9).Value = 3857.7354
$ws.Cells.Item(34, 10).Value = 28193.5
$ws.Cells.Item(34, 11).Value = 3857.7354
$ws.Cells.Item(34, 12).Value = 28193.5
$ws.Cells.Item(34, 13).Value = -3655.7354

$ws.Cells.Item(105, 8).Value = 62107.5
$ws.Cells.Item(105, 9).Value = 73329
$ws.Cells.Item(105, 10).Value = 6000
$ws.Cells.Item(105, 11).Value = 73329
$ws.Cells.Item(105, 12).Value = 6000
$ws.Cells.Item(105, 13).Value = -71582
$ws.Cells.Item(105, 14).Value = -9494

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(10, 8).Value = 273.73334
$ws.Cells.Item(10, 9).Value = 264.44446
$ws.Cells.Item(10, 10).Value = 287.66666
$ws.Cells.Item(10, 11).Value = 793.33338
$ws.Cells.Item(10, 12).Value = 862.9999799999999
$ws.Cells.Item(10, 13).Value = -654.33338

$ws.Cells.Item(17, 8).Value = 199.42857
$ws.Cells.Item(17, 9).Value = 38.8
$ws.Cells.Item(17, 10).Value = 601
$ws.Cells.Item(17, 11).Value = 116.4
$ws.Cells.Item(17, 12).Value = 1803
$ws.Cells.Item(17, 13).Value = 52.60000000000001
$ws.Cells.Item(17, 14).Value = -2141

$ws.Cells.Item(37, 8).Value = 97184.5
$ws.Cells.Item(37, 9).Value = 0
$ws.Cells.Item(37, 10).Value = 97184.5
$ws.Cells.Item(37, 11).Value = 0
$ws.Cells.Item(37, 12).Value = 291553.5
$ws.Cells.Item(37, 14).Value = -291777.5

$ws.Cells.Item(48, 8).Value = 7350
$ws.Cells.Item(48, 9).Value = 0
$ws.Cells.Item(48, 10).Value = 7350
$ws.Cells.Item(48, 11).Value = 0
$ws.Cells.Item(48, 12).Value = 22050
$ws.Cells.Item(48, 14).Value = -22550

$ws.Cells.Item(80, 8).Value = 5495
$ws.Cells.Item(80, 9).Value = 4995
$ws.Cells.Item(80, 10).Value = 5995
$ws.Cells.Item(80, 11).Value = 14985
$ws.Cells.Item(80, 12).Value = 17985
$ws.Cells.Item(80, 13).Value = -14049
$ws.Cells.Item(80, 14).Value = -19857

$ws.Cells.Item(83, 8).Value = 5495
$ws.Cells.Item(83, 9).Value = 4995
$ws.Cells.Item(83, 10).Value = 5995
$ws.Cells.Item(83, 11).Value = 44955
$ws.Cells.Item(83, 12).Value = 53955
$ws.Cells.Item(83, 13).Value = -40275
$ws.Cells.Item(83, 14).Value = -63315

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(47, 8).Value = 47500
$ws.Cells.Item(47, 9).Value = 47000
$ws.Cells.Item(47, 10).Value = 48000
$ws.Cells.Item(47, 11).Value = 47000
$ws.Cells.Item(47, 12).Value = 48000
$ws.Cells.Item(47, 13).Value = -46432
$ws.Cells.Item(47, 14).Value = -49136

$ws.Cells.Item(113, 8).Value = 5602.1304
$ws.Cells.Item(113, 9).Value = 4829.353
$ws.Cells.Item(113, 10).Value = 7791.6665
$ws.Cells.Item(113, 11).Value = 4829.353
$ws.Cells.Item(113, 12).Value = 7791.6665
$ws.Cells.Item(113, 13).Value = -2659.353

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 1771.7778
$ws.Cells.Item(68, 9).Value = 1618.25
$ws.Cells.Item(68, 10).Value = 3000
$ws.Cells.Item(68, 11).Value = 1618.25
$ws.Cells.Item(68, 12).Value = 3000
$ws.Cells.Item(68, 13).Value = -869.25
$ws.Cells.Item(68, 14).Value = -4498

$ws.Cells.Item(71, 8).Value = 1771.7778
$ws.Cells.Item(71, 9).Value = 1618.25
$ws.Cells.Item(71, 10).Value = 3000
$ws.Cells.Item(71, 11).Value = 8091.25
$ws.Cells.Item(71, 12).Value = 15000
$ws.Cells.Item(71, 13).Value = -4347.25
$ws.Cells.Item(71, 14).Value = -22488

$ws.Cells.Item(93, 8).Value = 1914.6154
$ws.Cells.Item(93, 9).Value = 1839
$ws.Cells.Item(93, 10).Value = 2166.6667
$ws.Cells.Item(93, 11).Value = 1839
$ws.Cells.Item(93, 12).Value = 2166.6667
$ws.Cells.Item(93, 13).Value = -591
$ws.Cells.Item(93, 14).Value = -4662.6667

$ws.Cells.Item(132, 8).Value = 11652.121
$ws.Cells.Item(132, 9).Value = 12846.962
$ws.Cells.Item(132, 10).Value = 7214.143
$ws.Cells.Item(132, 11).Value = 38540.886
$ws.Cells.Item(132, 12).Value = 21642.429
$ws.Cells.Item(132, 13).Value = -36010.886

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 0
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 13).Value = $null
$ws.Cells.Item(62, 14).Value = $null

$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).Value = $null
$ws.Cells.Item(65, 14).Value = $null

$ws.Cells.Item(126, 8).Value = 3024.8333
$ws.Cells.Item(126, 9).Value = 2936.2727
$ws.Cells.Item(126, 10).Value = 3999
$ws.Cells.Item(126, 11).Value = 8808.8181
$ws.Cells.Item(126, 12).Value = 11997
$ws.Cells.Item(126, 13).Value = -6338.8181

$ws.Cells.Item(136, 8).Value = 13919065
$ws.Cells.Item(136, 9).Value = 16701678
$ws.Cells.Item(136, 10).Value = 5999.25
$ws.Cells.Item(136, 11).Value = 50105034
$ws.Cells.Item(136, 12).Value = 17997.75
$ws.Cells.Item(136, 13).Value = -50102484
$ws.Cells.Item(136, 14).Value = -23097.75

$ws.Cells.Item(140, 8).Value = 0
$ws.Cells.Item(140, 9).Value = 0
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(140, 14).Value = $null

$ws.Cells.Item(141, 8).Value = 84000
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 84000
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 84000
$ws.Cells.Item(141, 14).Value = -94360
